$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1190.4  # H28: 1325.3334 -> 1190.4
$ws.Cells.Item(28, 9).Value = 1236.5  # I28: 1389.2 -> 1236.5
$ws.Cells.Item(28, 11).Value = 1236.5  # K28: 1389.2 -> 1236.5
$ws.Cells.Item(28, 13).Value = -751.5  # M28: -904.2 -> -751.5
$ws.Cells.Item(38, 8).Value = 6123.2  # H38: 3101.5715 -> 6123.2
$ws.Cells.Item(38, 9).Value = 872  # I38: 618.5 -> 872
$ws.Cells.Item(38, 10).Value = 14000  # J38: 18000 -> 14000
$ws.Cells.Item(38, 11).Value = 2616  # K38: 1855.5 -> 2616
$ws.Cells.Item(38, 12).Value = 42000  # L38: 54000 -> 42000
$ws.Cells.Item(38, 13).Value = -2244  # M38: -1483.5 -> -2244
$ws.Cells.Item(38, 14).Value = -42744  # N38: -54744 -> -42744
$ws.Cells.Item(103, 8).Value = 966.6667  # H103: 1000 -> 966.6667
$ws.Cells.Item(103, 9).Value = 750  # I103: 0 -> 750
$ws.Cells.Item(103, 10).Value = 993.75  # J103: 1000 -> 993.75
$ws.Cells.Item(103, 11).Value = 2250  # K103: 0 -> 2250
$ws.Cells.Item(103, 12).Value = 2981.25  # L103: 3000 -> 2981.25
$ws.Cells.Item(103, 13).Value = -1664  # M103: None -> -1664
$ws.Cells.Item(103, 14).Value = -4153.25  # N103: -4172 -> -4153.25
$ws.Cells.Item(111, 8).Value = 2515  # H111: 2182.0715 -> 2515
$ws.Cells.Item(111, 9).Value = 2515  # I111: 2280.6924 -> 2515
$ws.Cells.Item(111, 10).Value = 0  # J111: 900 -> 0
$ws.Cells.Item(111, 11).Value = 7545  # K111: 6842.0772 -> 7545
$ws.Cells.Item(111, 12).Value = 0  # L111: 2700 -> 0
$ws.Cells.Item(111, 13).Value = -4478  # M111: -3775.0772 -> -4478
$ws.Cells.Item(111, 14).ClearContents()  # N111: -8834 -> (removed)
$ws.Cells.Item(138, 8).Value = 2442.7576  # H138: 2442.7273 -> 2442.7576
$ws.Cells.Item(138, 10).Value = 3234.318  # J138: 3234.2727 -> 3234.318
$ws.Cells.Item(138, 12).Value = 9702.954000000002  # L138: 9702.8181 -> 9702.954000000002
$ws.Cells.Item(138, 14).Value = -19982.954  # N138: -19982.8181 -> -19982.954

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1631.9286  # H2: 1662.3334 -> 1631.9286
$ws.Cells.Item(2, 9).Value = 1680.5385  # I2: 1662.3334 -> 1680.5385
$ws.Cells.Item(2, 10).Value = 1000  # J2: 0 -> 1000
$ws.Cells.Item(2, 11).Value = 1680.5385  # K2: 1662.3334 -> 1680.5385
$ws.Cells.Item(2, 12).Value = 1000  # L2: 0 -> 1000
$ws.Cells.Item(2, 13).Value = -1567.5385  # M2: -1549.3334 -> -1567.5385
$ws.Cells.Item(2, 14).Value = -1226  # N2: None -> -1226
$ws.Cells.Item(6, 8).Value = 5000  # H6: 0 -> 5000
$ws.Cells.Item(6, 10).Value = 5000  # J6: 0 -> 5000
$ws.Cells.Item(6, 12).Value = 5000  # L6: 0 -> 5000
$ws.Cells.Item(6, 14).Value = -5346  # N6: None -> -5346
$ws.Cells.Item(51, 8).Value = 38495  # H51: 0 -> 38495
$ws.Cells.Item(51, 10).Value = 38495  # J51: 0 -> 38495
$ws.Cells.Item(51, 12).Value = 38495  # L51: 0 -> 38495
$ws.Cells.Item(51, 14).Value = -40007  # N51: None -> -40007
$ws.Cells.Item(81, 8).Value = 0  # H81: 122000 -> 0
$ws.Cells.Item(81, 10).Value = 0  # J81: 122000 -> 0
$ws.Cells.Item(81, 12).Value = 0  # L81: 122000 -> 0
$ws.Cells.Item(81, 14).ClearContents()  # N81: -123996 -> (removed)
$ws.Cells.Item(84, 8).Value = 0  # H84: 122000 -> 0
$ws.Cells.Item(84, 10).Value = 0  # J84: 122000 -> 0
$ws.Cells.Item(84, 12).Value = 0  # L84: 366000 -> 0
$ws.Cells.Item(84, 14).ClearContents()  # N84: -375984 -> (removed)
$ws.Cells.Item(113, 8).Value = 93299.664  # H113: 69977.8 -> 93299.664
$ws.Cells.Item(113, 10).Value = 93299.664  # J113: 69977.8 -> 93299.664
$ws.Cells.Item(113, 12).Value = 93299.664  # L113: 69977.8 -> 93299.664
$ws.Cells.Item(113, 14).Value = -101977.664  # N113: -78655.8 -> -101977.664
$ws.Cells.Item(116, 8).Value = 1631.9286  # H116: 1662.3334 -> 1631.9286
$ws.Cells.Item(116, 9).Value = 1680.5385  # I116: 1662.3334 -> 1680.5385
$ws.Cells.Item(116, 10).Value = 1000  # J116: 0 -> 1000
$ws.Cells.Item(116, 11).Value = 1680.5385  # K116: 1662.3334 -> 1680.5385
$ws.Cells.Item(116, 12).Value = 1000  # L116: 0 -> 1000
$ws.Cells.Item(116, 13).Value = 613.4614999999999  # M116: 631.6666 -> 613.4614999999999
$ws.Cells.Item(116, 14).Value = -5588  # N116: None -> -5588

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1631.9286  # H3: 1662.3334 -> 1631.9286
$ws.Cells.Item(3, 9).Value = 1680.5385  # I3: 1662.3334 -> 1680.5385
$ws.Cells.Item(3, 10).Value = 1000  # J3: 0 -> 1000
$ws.Cells.Item(3, 11).Value = 1680.5385  # K3: 1662.3334 -> 1680.5385
$ws.Cells.Item(3, 12).Value = 1000  # L3: 0 -> 1000
$ws.Cells.Item(3, 13).Value = -1566.5385  # M3: -1548.3334 -> -1566.5385
$ws.Cells.Item(3, 14).Value = -1228  # N3: None -> -1228
$ws.Cells.Item(107, 8).Value = 2213.5  # H107: 2131 -> 2213.5
$ws.Cells.Item(107, 9).Value = 2118  # I107: 2044.2858 -> 2118
$ws.Cells.Item(107, 10).Value = 2500  # J107: 2333.3333 -> 2500
$ws.Cells.Item(107, 11).Value = 2118  # K107: 2044.2858 -> 2118
$ws.Cells.Item(107, 12).Value = 2500  # L107: 2333.3333 -> 2500
$ws.Cells.Item(107, 13).Value = -198  # M107: -124.2858000000001 -> -198
$ws.Cells.Item(107, 14).Value = -6340  # N107: -6173.3333 -> -6340
$ws.Cells.Item(132, 8).Value = 70000  # H132: 69000 -> 70000
$ws.Cells.Item(132, 10).Value = 70000  # J132: 69000 -> 70000
$ws.Cells.Item(132, 12).Value = 70000  # L132: 69000 -> 70000
$ws.Cells.Item(132, 14).Value = -80120  # N132: -79120 -> -80120

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 571244.9399999999  # H31: 557681.8 -> 571244.9399999999
$ws.Cells.Item(31, 9).Value = 1855.3158  # I31: 1842.25 -> 1855.3158
$ws.Cells.Item(31, 11).Value = 1855.3158  # K31: 1842.25 -> 1855.3158
$ws.Cells.Item(31, 13).Value = -1560.3158  # M31: -1547.25 -> -1560.3158
$ws.Cells.Item(32, 8).Value = 4000  # H32: 0 -> 4000
$ws.Cells.Item(32, 9).Value = 4000  # I32: 0 -> 4000
$ws.Cells.Item(32, 11).Value = 4000  # K32: 0 -> 4000
$ws.Cells.Item(32, 13).Value = -3684  # M32: None -> -3684
$ws.Cells.Item(34, 8).Value = 571244.9399999999  # H34: 557681.8 -> 571244.9399999999
$ws.Cells.Item(34, 9).Value = 1855.3158  # I34: 1842.25 -> 1855.3158
$ws.Cells.Item(34, 11).Value = 1855.3158  # K34: 1842.25 -> 1855.3158
$ws.Cells.Item(34, 13).Value = -1653.3158  # M34: -1640.25 -> -1653.3158
$ws.Cells.Item(69, 8).Value = 86894.5  # H69: 54238.2 -> 86894.5
$ws.Cells.Item(69, 10).Value = 86894.5  # J69: 54238.2 -> 86894.5
$ws.Cells.Item(69, 12).Value = 86894.5  # L69: 54238.2 -> 86894.5
$ws.Cells.Item(69, 14).Value = -88392.5  # N69: -55736.2 -> -88392.5
$ws.Cells.Item(72, 8).Value = 86894.5  # H72: 54238.2 -> 86894.5
$ws.Cells.Item(72, 10).Value = 86894.5  # J72: 54238.2 -> 86894.5
$ws.Cells.Item(72, 12).Value = 260683.5  # L72: 162714.6 -> 260683.5
$ws.Cells.Item(72, 14).Value = -268171.5  # N72: -170202.6 -> -268171.5
$ws.Cells.Item(87, 8).Value = 65000  # H87: 81216.336 -> 65000
$ws.Cells.Item(87, 10).Value = 0  # J87: 89324.5 -> 0
$ws.Cells.Item(87, 12).Value = 0  # L87: 89324.5 -> 0
$ws.Cells.Item(87, 14).ClearContents()  # N87: -91696.5 -> (removed)
$ws.Cells.Item(90, 8).Value = 65000  # H90: 81216.336 -> 65000
$ws.Cells.Item(90, 10).Value = 0  # J90: 89324.5 -> 0
$ws.Cells.Item(90, 12).Value = 0  # L90: 267973.5 -> 0
$ws.Cells.Item(90, 14).ClearContents()  # N90: -279829.5 -> (removed)
$ws.Cells.Item(107, 8).Value = 2843.8333  # H107: 2922.6 -> 2843.8333
$ws.Cells.Item(107, 9).Value = 2810  # I107: 2900 -> 2810
$ws.Cells.Item(107, 11).Value = 2810  # K107: 2900 -> 2810
$ws.Cells.Item(107, 13).Value = -890  # M107: -980 -> -890
$ws.Cells.Item(127, 8).Value = 84995  # H127: 119000 -> 84995
$ws.Cells.Item(127, 10).Value = 84995  # J127: 119000 -> 84995
$ws.Cells.Item(127, 12).Value = 84995  # L127: 119000 -> 84995
$ws.Cells.Item(127, 14).Value = -94915  # N127: -128920 -> -94915

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 5136437  # H4: 4540069.5 -> 5136437
$ws.Cells.Item(4, 9).Value = 5333410.5  # I4: 4869640.5 -> 5333410.5
$ws.Cells.Item(4, 10).Value = 1000000  # J4: 750000 -> 1000000
$ws.Cells.Item(4, 11).Value = 16000231.5  # K4: 14608921.5 -> 16000231.5
$ws.Cells.Item(4, 12).Value = 3000000  # L4: 2250000 -> 3000000
$ws.Cells.Item(4, 13).Value = -16000119.5  # M4: -14608809.5 -> -16000119.5
$ws.Cells.Item(4, 14).Value = -3000224  # N4: -2250224 -> -3000224
$ws.Cells.Item(38, 8).Value = 58.8125  # H38: 48.61905 -> 58.8125
$ws.Cells.Item(38, 9).Value = 11  # I38: 17.8 -> 11
$ws.Cells.Item(38, 10).Value = 62  # J38: 58.25 -> 62
$ws.Cells.Item(38, 11).Value = 33  # K38: 53.40000000000001 -> 33
$ws.Cells.Item(38, 12).Value = 186  # L38: 174.75 -> 186
$ws.Cells.Item(38, 13).Value = 314  # M38: 293.6 -> 314
$ws.Cells.Item(38, 14).Value = -880  # N38: -868.75 -> -880
$ws.Cells.Item(114, 8).Value = 1689.7  # H114: 1808.909 -> 1689.7
$ws.Cells.Item(114, 10).Value = 0  # J114: 3001 -> 0
$ws.Cells.Item(114, 12).Value = 0  # L114: 9003 -> 0
$ws.Cells.Item(114, 14).ClearContents()  # N114: -15511 -> (removed)

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 10570.863  # H5: 10679.909 -> 10570.863
$ws.Cells.Item(5, 9).Value = 7502.579  # I5: 7628.8423 -> 7502.579
$ws.Cells.Item(5, 11).Value = 7502.579  # K5: 7628.8423 -> 7502.579
$ws.Cells.Item(5, 13).Value = -7390.579  # M5: -7516.8423 -> -7390.579
$ws.Cells.Item(19, 8).Value = 0  # H19: 5 -> 0
$ws.Cells.Item(19, 9).Value = 0  # I19: 5 -> 0
$ws.Cells.Item(19, 11).Value = 0  # K19: 5 -> 0
$ws.Cells.Item(19, 13).ClearContents()  # M19: 283 -> (removed)
$ws.Cells.Item(24, 8).Value = 33898.332  # H24: 23335.666 -> 33898.332
$ws.Cells.Item(24, 9).Value = 14503  # I24: 0 -> 14503
$ws.Cells.Item(24, 10).Value = 37777.4  # J24: 23335.666 -> 37777.4
$ws.Cells.Item(24, 11).Value = 14503  # K24: 0 -> 14503
$ws.Cells.Item(24, 12).Value = 37777.4  # L24: 23335.666 -> 37777.4
$ws.Cells.Item(24, 13).Value = -14330  # M24: None -> -14330
$ws.Cells.Item(24, 14).Value = -38123.4  # N24: -23681.666 -> -38123.4
$ws.Cells.Item(62, 8).Value = 99957.664  # H62: 100437 -> 99957.664
$ws.Cells.Item(62, 10).Value = 99957.664  # J62: 100437 -> 99957.664
$ws.Cells.Item(62, 12).Value = 99957.664  # L62: 100437 -> 99957.664
$ws.Cells.Item(62, 14).Value = -101329.664  # N62: -101809 -> -101329.664
$ws.Cells.Item(65, 8).Value = 99957.664  # H65: 100437 -> 99957.664
$ws.Cells.Item(65, 10).Value = 99957.664  # J65: 100437 -> 99957.664
$ws.Cells.Item(65, 12).Value = 299872.992  # L65: 301311 -> 299872.992
$ws.Cells.Item(65, 14).Value = -306736.992  # N65: -308175 -> -306736.992
$ws.Cells.Item(68, 8).Value = 0  # H68: 20000 -> 0
$ws.Cells.Item(68, 9).Value = 0  # I68: 20000 -> 0
$ws.Cells.Item(68, 11).Value = 0  # K68: 20000 -> 0
$ws.Cells.Item(68, 13).ClearContents()  # M68: -19189 -> (removed)
$ws.Cells.Item(71, 8).Value = 0  # H71: 20000 -> 0
$ws.Cells.Item(71, 9).Value = 0  # I71: 20000 -> 0
$ws.Cells.Item(71, 11).Value = 0  # K71: 60000 -> 0
$ws.Cells.Item(71, 13).ClearContents()  # M71: -55944 -> (removed)
$ws.Cells.Item(74, 8).Value = 0  # H74: 35000 -> 0
$ws.Cells.Item(74, 10).Value = 0  # J74: 35000 -> 0
$ws.Cells.Item(74, 12).Value = 0  # L74: 35000 -> 0
$ws.Cells.Item(74, 14).ClearContents()  # N74: -36872 -> (removed)
$ws.Cells.Item(77, 8).Value = 0  # H77: 35000 -> 0
$ws.Cells.Item(77, 10).Value = 0  # J77: 35000 -> 0
$ws.Cells.Item(77, 12).Value = 0  # L77: 105000 -> 0
$ws.Cells.Item(77, 14).ClearContents()  # N77: -114360 -> (removed)
$ws.Cells.Item(100, 8).Value = 59995  # H100: 115000 -> 59995
$ws.Cells.Item(100, 10).Value = 59995  # J100: 115000 -> 59995
$ws.Cells.Item(100, 12).Value = 59995  # L100: 115000 -> 59995
$ws.Cells.Item(100, 14).Value = -62159  # N100: -117164 -> -62159

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2963.4194  # H46: 2966.9678 -> 2963.4194
$ws.Cells.Item(46, 9).Value = 2872.4783  # I46: 2799.0417 -> 2872.4783
$ws.Cells.Item(46, 10).Value = 3224.875  # J46: 3542.7144 -> 3224.875
$ws.Cells.Item(46, 11).Value = 2872.4783  # K46: 2799.0417 -> 2872.4783
$ws.Cells.Item(46, 12).Value = 3224.875  # L46: 3542.7144 -> 3224.875
$ws.Cells.Item(46, 13).Value = -2684.4783  # M46: -2611.0417 -> -2684.4783
$ws.Cells.Item(46, 14).Value = -3600.875  # N46: -3918.7144 -> -3600.875
$ws.Cells.Item(48, 8).Value = 37495  # H48: 0 -> 37495
$ws.Cells.Item(48, 10).Value = 37495  # J48: 0 -> 37495
$ws.Cells.Item(48, 12).Value = 37495  # L48: 0 -> 37495
$ws.Cells.Item(48, 14).Value = -38817  # N48: None -> -38817
$ws.Cells.Item(63, 8).Value = 82000  # H63: 0 -> 82000
$ws.Cells.Item(63, 10).Value = 82000  # J63: 0 -> 82000
$ws.Cells.Item(63, 12).Value = 82000  # L63: 0 -> 82000
$ws.Cells.Item(63, 14).Value = -83498  # N63: None -> -83498
$ws.Cells.Item(66, 8).Value = 82000  # H66: 0 -> 82000
$ws.Cells.Item(66, 10).Value = 82000  # J66: 0 -> 82000
$ws.Cells.Item(66, 12).Value = 246000  # L66: 0 -> 246000
$ws.Cells.Item(66, 14).Value = -253488  # N66: None -> -253488
$ws.Cells.Item(74, 8).Value = 64994  # H74: 118000 -> 64994
$ws.Cells.Item(74, 10).Value = 64994  # J74: 118000 -> 64994
$ws.Cells.Item(74, 12).Value = 64994  # L74: 118000 -> 64994
$ws.Cells.Item(74, 14).Value = -66990  # N74: -119996 -> -66990
$ws.Cells.Item(77, 8).Value = 64994  # H77: 118000 -> 64994
$ws.Cells.Item(77, 10).Value = 64994  # J77: 118000 -> 64994
$ws.Cells.Item(77, 12).Value = 194982  # L77: 354000 -> 194982
$ws.Cells.Item(77, 14).Value = -204966  # N77: -363984 -> -204966
$ws.Cells.Item(132, 8).Value = 224901.67  # H132: 121528.06 -> 224901.67
$ws.Cells.Item(132, 9).Value = 169018.33  # I132: 202222.4 -> 169018.33
$ws.Cells.Item(132, 10).Value = 336668.34  # J132: 87905.414 -> 336668.34
$ws.Cells.Item(132, 11).Value = 507054.99  # K132: 606667.2 -> 507054.99
$ws.Cells.Item(132, 12).Value = 1010005.02  # L132: 263716.242 -> 1010005.02
$ws.Cells.Item(132, 13).Value = -504524.99  # M132: -604137.2 -> -504524.99
$ws.Cells.Item(132, 14).Value = -1015065.02  # N132: -268776.242 -> -1015065.02
$ws.Cells.Item(136, 8).Value = 47185.703  # H136: 50633.24 -> 47185.703
$ws.Cells.Item(136, 9).Value = 4414.5454  # I136: 4537.6 -> 4414.5454
$ws.Cells.Item(136, 10).Value = 76590.875  # J136: 81363.664 -> 76590.875
$ws.Cells.Item(136, 11).Value = 13243.6362  # K136: 13612.8 -> 13243.6362
$ws.Cells.Item(136, 12).Value = 229772.625  # L136: 244090.992 -> 229772.625
$ws.Cells.Item(136, 13).Value = -10693.6362  # M136: -11062.8 -> -10693.6362
$ws.Cells.Item(136, 14).Value = -234872.625  # N136: -249190.992 -> -234872.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 1000  # H9: 503 -> 1000
$ws.Cells.Item(9, 9).Value = 0  # I9: 6 -> 0
$ws.Cells.Item(9, 11).Value = 0  # K9: 6 -> 0
$ws.Cells.Item(9, 13).ClearContents()  # M9: 134 -> (removed)
$ws.Cells.Item(11, 8).Value = 972.5  # H11: 4799 -> 972.5
$ws.Cells.Item(11, 9).Value = 995  # I11: 0 -> 995
$ws.Cells.Item(11, 10).Value = 950  # J11: 4799 -> 950
$ws.Cells.Item(11, 11).Value = 995  # K11: 0 -> 995
$ws.Cells.Item(11, 12).Value = 950  # L11: 4799 -> 950
$ws.Cells.Item(11, 13).Value = -853  # M11: None -> -853
$ws.Cells.Item(11, 14).Value = -1234  # N11: -5083 -> -1234
$ws.Cells.Item(107, 8).Value = 16130311  # H107: 16667978 -> 16130311
$ws.Cells.Item(107, 9).Value = 20834656  # I107: 21740498 -> 20834656
$ws.Cells.Item(107, 11).Value = 62503968  # K107: 65221494 -> 62503968
$ws.Cells.Item(107, 13).Value = -62502048  # M107: -65219574 -> -62502048
$ws.Cells.Item(132, 8).Value = 1951.0435  # H132: 2070.3333 -> 1951.0435
$ws.Cells.Item(132, 9).Value = 2021.6818  # I132: 2080.8572 -> 2021.6818
$ws.Cells.Item(132, 10).Value = 397  # J132: 1996.6666 -> 397
$ws.Cells.Item(132, 11).Value = 6065.0454  # K132: 6242.571599999999 -> 6065.0454
$ws.Cells.Item(132, 12).Value = 1191  # L132: 5989.9998 -> 1191
$ws.Cells.Item(132, 13).Value = -3535.0454  # M132: -3712.571599999999 -> -3535.0454
$ws.Cells.Item(136, 8).Value = 1852.8  # H136: 2191.6667 -> 1852.8
$ws.Cells.Item(136, 9).Value = 1066  # I136: 1255.5555 -> 1066
$ws.Cells.Item(136, 11).Value = 3198  # K136: 3766.6665 -> 3198
$ws.Cells.Item(136, 13).Value = -648  # M136: -1216.6665 -> -648
